$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7264
$ws1.Range("F5").Value = 267
$ws1.Range("F6").Value = 424
$ws1.Range("F7").Value = 3748

# Sheet "全部类型" (all types) — same events repeated, update the same values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7264
$ws4.Range("F7").Value = 267
$ws4.Range("F8").Value = 424
$ws4.Range("F9").Value = 3748
